$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with new values
$ws.Range("A1").Value = "Copper Horse - Warhorse 5lb"
$ws.Range("B1").Value = 18

$ws.Range("A2").Value = "Copper Horse - Warhorse Blend (12oz)"
$ws.Range("B2").Value = 10

$ws.Range("A3").Value = "Copper Horse - Clocktower Espresso (12oz)"
$ws.Range("B3").Value = 7

$ws.Range("A4").Value = "Copper Horse - Carriage House Blend (12oz)"
$ws.Range("B4").Value = 7

$ws.Range("A1").Select()
